$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data to the endpoint table
$ws.Range("A17").Value = "Refresh Token"
$ws.Range("B17").Value = "/api/v1/get-access-token"
$ws.Range("C17").Value = '{"refresh_token":"string"}'
$ws.Range("D17").Value = '{"token":"string"}'

# Match formatting of the surrounding table (border + left/vcenter alignment like column A above)
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("A17").HorizontalAlignment = -4131
$ws.Range("A17").VerticalAlignment = -4108

$ws.Range("B17").Borders.LineStyle = 1

$ws.Range("C17").Borders.LineStyle = 1
$ws.Range("C17").WrapText = $true

$ws.Range("D17").Borders.LineStyle = 1

# Update the view: scroll so row 13 is at the top, and select C5
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("C5").Select()
